$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: update VAT ID cell (E2) text (removes stray _x000D_ CR markers that
# came from a pasted Outlook body, keeping a single trailing line break) and
# re-flow row height for the shorter text.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "DE1567890`n"
$ws.Rows.Item(2).RowHeight = 45

# ---------------------------------------------------------------------------
# Row 3: new email, duplicate of row 2's data/format (same sender, company,
# address, VAT id, email & subject, different timestamp).
# ---------------------------------------------------------------------------
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)
$ws.Range("A3").Value = 45854.104166666664
$ws.Range("B3").Value = "Moris Mwai"
$ws.Range("C3").Value = "Tech-Neo GmbH"
$ws.Range("D3").Value = "Am main City, Germany"
$ws.Range("E3").Value = "DE1567890`n"
$ws.Range("F3").Value = "morismwai1@gmail.com"
$ws.Range("G3").Value = "Partnership Offer"
$ws.Rows.Item(3).RowHeight = 45

# ---------------------------------------------------------------------------
# Rows 5-13: additional processed emails appended by the dispatcher (note -
# row 4 is intentionally skipped, matching the source queue export).
# ---------------------------------------------------------------------------
$dates = @(45854.23364583333, 45854.289583333331, 45854.293391203704, 45854.307812500003, 45854.317210648151, 45854.348321759258, 45854.356712962966, 45854.362986111111, 45854.369895833333)
$rows = @(5, 6, 7, 8, 9, 10, 11, 12, 13)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Range("A2").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("A$r").Value = $dates[$i]
    $ws.Range("B$r").Value = "Moris Mwai"
    $ws.Range("C$r").Value = "Tech-Neo GmbH"
    $ws.Range("D$r").Value = "Am main City, Germany"
    $ws.Range("E$r").Value = "DE1567890"
    $ws.Range("F$r").Value = "morismwai1@gmail.com"
    $ws.Range("G$r").Value = "Partnership Offer"
}

# ---------------------------------------------------------------------------
# Selection / active cell bookkeeping to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("E2").Select() | Out-Null
